# Updated cryptos list on Mon May 20 14:45:48 UTC 2024 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns for the crypto ranking
# sheet, and reflects Arweave/Filecoin swapping ranking positions (rows 37/38).
#
# D-column values that parse as plain numbers (e.g. "175.60") get
# apostrophe-prefixed so Excel stores them as text (matching the sheet's
# existing text-formatted price strings, e.g. "66.850.68" keeps its
# thousands-dot formatting instead of being coerced to a float), then the
# quote-prefix style is reset back to Normal so no stray per-cell formatting
# is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.850.68'
$ws.Range('E2').Value = '  -0.22%  '
$ws.Range('D3').Value = '3.082.20'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = "'570.31"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.10%  '
$ws.Range('D6').Value = "'175.60"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.49%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '3.081.85'
$ws.Range('E8').Value = '  +0.26%  '
$ws.Range('E9').Value = '  -0.24%  '
$ws.Range('E10').Value = '  +0.46%  '
$ws.Range('E11').Value = '  +0.52%  '
$ws.Range('E12').Value = '  -1.24%  '
$ws.Range('E13').Value = '  -0.82%  '
$ws.Range('D14').Value = "'35.79"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.72%  '
$ws.Range('E15').Value = '  +0.82%  '
$ws.Range('D16').Value = '3.595.50'
$ws.Range('E16').Value = '  +0.14%  '
$ws.Range('D17').Value = '66.783.88'
$ws.Range('E17').Value = '  -0.25%  '
$ws.Range('E18').Value = '  -0.55%  '
$ws.Range('D19').Value = '3.083.32'
$ws.Range('E19').Value = '  +0.13%  '
$ws.Range('D20').Value = "'16.46"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.06%  '
$ws.Range('D21').Value = "'483.29"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.63%  '
$ws.Range('D22').Value = "'7.68"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.15%  '
$ws.Range('E23').Value = '  -0.61%  '
$ws.Range('D24').Value = "'83.24"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.55%  '
$ws.Range('D25').Value = "'12.68"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.36%  '
$ws.Range('D26').Value = "'2.22"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.47%  '
$ws.Range('E27').Value = '  +0.93%  '
$ws.Range('E28').Value = '  -0.06%  '
$ws.Range('D29').Value = "'7.82"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.58%  '
$ws.Range('E30').Value = '  -0.01%  '
$ws.Range('E31').Value = '  -1.62%  '
$ws.Range('D32').Value = "'27.87"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.19%  '
$ws.Range('E33').Value = '  +0.08%  '
$ws.Range('D34').Value = '0.0₃0921'
$ws.Range('E34').Value = '  +1.83%  '
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('D36').Value = "'0.945"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.00%  '
$ws.Range('B37').Value = 'Arweave'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D37').Value = "'46.73"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.17%  '
$ws.Range('B38').Value = 'Filecoin'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D38').Value = "'5.52"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.83%  '
$ws.Range('E39').Value = '  +2.87%  '
$ws.Range('D40').Value = "'48.91"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.45%  '
$ws.Range('E41').Value = '  +0.22%  '
$ws.Range('E42').Value = '  -0.37%  '
$ws.Range('D43').Value = "'8.20"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.33%  '
$ws.Range('E44').Value = '  +8.52%  '
$ws.Range('D45').Value = '2.792.52'
$ws.Range('E45').Value = '  +0.64%  '
$ws.Range('D46').Value = "'368.86"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.05%  '
$ws.Range('D47').Value = "'134.81"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.66%  '
$ws.Range('D48').Value = "'0.0341"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.83%  '
$ws.Range('E49').Value = '  -0.01%  '
$ws.Range('D50').Value = "'25.64"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.25%  '
$ws.Range('D51').Value = "'2.28"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +6.04%  '
